$d = $word.ActiveDocument

# 1) Change the "Introduction" heading paragraph's style to "Body Text First Indent"
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Introduction*") {
        $p.Style = "Body Text First Indent"
    }
}
